$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so that numeric-looking
# strings like "556.85" are not auto-converted to numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '66.560.68'
$ws.Range("E2").Value = '  -5.13%  '

$ws.Range("D3").Value = '3.357.14'
$ws.Range("E3").Value = '  -5.99%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '556.85'
$ws.Range("E5").Value = '  -5.32%  '

$ws.Range("D6").Value = '183.28'
$ws.Range("E6").Value = '  -7.93%  '

$ws.Range("D7").Value = '0.598'
$ws.Range("E7").Value = '  -4.07%  '

$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").Value = '3.347.53'
$ws.Range("E9").Value = '  -5.78%  '

$ws.Range("D10").Value = '0.188'
$ws.Range("E10").Value = '  -12.10%  '

$ws.Range("D11").Value = '0.593'
$ws.Range("E11").Value = '  -7.25%  '

$ws.Range("D12").Value = '47.80'
$ws.Range("E12").Value = '  -9.71%  '

$ws.Range("E13").Value = '  -10.33%  '

$ws.Range("D14").Value = '8.70'
$ws.Range("E14").Value = '  -8.61%  '

$ws.Range("D15").Value = '3.883.04'
$ws.Range("E15").Value = '  -6.41%  '

$ws.Range("D16").Value = '601.13'
$ws.Range("E16").Value = '  -13.48%  '

$ws.Range("D17").Value = '66.341.02'
$ws.Range("E17").Value = '  -5.56%  '

$ws.Range("D18").Value = '3.344.92'
$ws.Range("E18").Value = '  -6.86%  '

$ws.Range("E19").Value = '  -4.49%  '

$ws.Range("D20").Value = '17.85'
$ws.Range("E20").Value = '  -5.72%  '

$ws.Range("D21").Value = '11.66'
$ws.Range("E21").Value = '  -7.83%  '

$ws.Range("D22").Value = '0.910'
$ws.Range("E22").Value = '  -7.83%  '

$ws.Range("D23").Value = '16.83'
$ws.Range("E23").Value = '  -7.51%  '

$ws.Range("D24").Value = '5.04'
$ws.Range("E24").Value = '  -4.63%  '

$ws.Range("D25").Value = '97.93'
$ws.Range("E25").Value = '  -11.69%  '

$ws.Range("D26").Value = '4.06'
$ws.Range("E26").Value = '  -9.23%  '

$ws.Range("D27").Value = '2.75'
$ws.Range("E27").Value = '  -7.87%  '

$ws.Range("D28").Value = '9.44'
$ws.Range("E28").Value = '  -9.62%  '

$ws.Range("D29").Value = '8.80'
$ws.Range("E29").Value = '  -11.19%  '

$ws.Range("D30").Value = '30.71'
$ws.Range("E30").Value = '  -10.45%  '

$ws.Range("D31").Value = '6.34'
$ws.Range("E31").Value = '  -9.76%  '

$ws.Range("D32").Value = '3.85'
$ws.Range("E32").Value = '  -13.06%  '

$ws.Range("D33").Value = '11.16'
$ws.Range("E33").Value = '  -8.16%  '

$ws.Range("E34").Value = '  -7.18%  '

$ws.Range("D35").Value = '3.831.15'
$ws.Range("E35").Value = '  +0.77%  '

$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").Value = '57.85'
$ws.Range("E36").Value = '  -8.76%  '

$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").Value = '532.08'
$ws.Range("E37").Value = '  +5.08%  '

$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.05%  '

$ws.Range("D39").Value = '3.66'
$ws.Range("E39").Value = '  +36.20%  '

$ws.Range("D40").Value = '3.39'
$ws.Range("E40").Value = '  -7.33%  '

$ws.Range("D41").Value = '0.0₃0724'
$ws.Range("E41").Value = '  -13.66%  '

$ws.Range("D42").Value = '2.71'
$ws.Range("E42").Value = '  -8.81%  '

$ws.Range("D43").Value = '0.127'
$ws.Range("E43").Value = '  -7.22%  '

$ws.Range("D44").Value = '0.349'
$ws.Range("E44").Value = '  -7.95%  '

$ws.Range("D45").Value = '32.56'
$ws.Range("E45").Value = '  -9.75%  '

$ws.Range("D46").Value = '0.0417'
$ws.Range("E46").Value = '  -10.60%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.16'
$ws.Range("E47").Value = '  -7.97%  '

$ws.Range("B48").Value = 'ThetaToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D48").Value = '2.68'
$ws.Range("E48").Value = '  -11.49%  '

$ws.Range("E49").Value = '  -7.08%  '

$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  -0.40%  '

$ws.Range("D51").Value = '7.69'
$ws.Range("E51").Value = '  -10.33%  '

# Restore default (unstyled) formatting on column D so the output cells
# do not retain an explicit number format style, matching the original file.
$dRange.ClearFormats()
